# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" sheet and the combined "全部类型" sheet, matching the refreshed
# data output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - rows 2-7
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 3457
$wsExhibit.Range("F3").Value = 30
$wsExhibit.Range("F4").Value = 72
$wsExhibit.Range("F5").Value = 1874
$wsExhibit.Range("F6").Value = 128
$wsExhibit.Range("F7").Value = 344

# Sheet "全部类型" - rows 2-6 and 8 (row 7 is an unrelated event, unchanged)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3457
$wsAll.Range("F3").Value = 30
$wsAll.Range("F4").Value = 72
$wsAll.Range("F5").Value = 1874
$wsAll.Range("F6").Value = 128
$wsAll.Range("F8").Value = 344
